# Append the next day's GSC export row to the "Chart" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Find the first empty row after the existing data (row 78 is currently
# the last data row, so this lands on row 79).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Write the date as literal text (matching the existing column A cells,
# which are plain strings rather than real date values), then drop the
# temporary text format so the cell keeps the sheet's default style.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-12-23"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 32
